$d = $word.ActiveDocument

# 1. Remove the existing (misplaced) _GoBack bookmark, wherever Word last left it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Rename the "Approach" heading to "Methods".
$d.Content.Find.Execute("Approach", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Methods", 2) | Out-Null

# 3. Re-create the _GoBack bookmark immediately after the new heading text,
#    inside the heading paragraph (this is where Word leaves it after the
#    user's last edit point). Inserting/removing a one-character placeholder
#    avoids the runtime's edge-case handling of zero-length ranges sitting
#    exactly on a paragraph boundary.
$headingRange = $d.Content
$headingRange.Find.Execute("Methods", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0) | Out-Null
$headingRange.Collapse(0) ;# wdCollapseEnd

$placeholderStart = $headingRange.Start
$headingRange.InsertAfter("X")

$placeholderRange = $d.Range($placeholderStart, $placeholderStart + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)

$placeholderRange2 = $d.Range($placeholderStart, $placeholderStart + 1)
$placeholderRange2.Delete()
